$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.977.79"
$ws.Range("E2").Value = "  +3.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.471.89"
$ws.Range("E3").Value = "  +5.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.13"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.56"
$ws.Range("E6").Value = "  +8.53%  "
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.469.13"
$ws.Range("E9").Value = "  +5.68%  "
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.69"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.68"
$ws.Range("E14").Value = "  +11.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.911.59"
$ws.Range("E15").Value = "  +5.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.832.10"
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.467.82"
$ws.Range("E18").Value = "  +5.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.54"
$ws.Range("E20").Value = "  +8.16%  "
$ws.Range("E21").Value = "  +3.64%  "
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.61"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.50"
$ws.Range("E27").Value = "  +6.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.39"
$ws.Range("E29").Value = "  +9.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.86"
$ws.Range("E30").Value = "  +12.35%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.84"
$ws.Range("E31").Value = "  +6.03%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0802"
$ws.Range("E32").Value = "  +9.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.15"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("E34").Value = "  +10.74%  "
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.79"
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "369.02"
$ws.Range("E37").Value = "  +10.67%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  +6.01%  "
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +10.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.61"
$ws.Range("E42").Value = "  +6.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "149.71"
$ws.Range("E43").Value = "  +7.75%  "
$ws.Range("E44").Value = "  +5.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.51"
$ws.Range("E45").Value = "  +6.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.600"
$ws.Range("E46").Value = "  +5.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0959"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"
$ws.Range("E49").Value = "  +4.35%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0232"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.96"
$ws.Range("E51").Value = "  +4.67%  "
